$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 37.63904266666666
$ws.Range("H2").Value = 112.917128
$ws.Range("I2").Value = 0.4850220755088102
$ws.Range("J2").Value = 0.4850220755088102
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.22896466666667
$ws.Range("N2").Value = 93.686894
$ws.Range("O2").Value = 0.2877106972998646
$ws.Range("P2").Value = 0.2877106972998646
$ws.Range("Q2").Value = 1175.428333524492
$ws.Range("R2").Value = 10578.85500172043
$ws.Range("S2").Value = 0.1395460395504673
$ws.Range("T2").Value = 0.1395460395504674
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 37.63904266666666
$ws.Range("H3").Value = 112.917128
$ws.Range("I3").Value = 0.4850220755088102
$ws.Range("J3").Value = 0.4850220755088102
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.44578266666667
$ws.Range("N3").Value = 121.337348
$ws.Range("O3").Value = 0.3726247238124506
$ws.Range("P3").Value = 0.3726247238124505
$ws.Range("Q3").Value = 1522.340539477394
$ws.Range("R3").Value = 13701.06485529654
$ws.Range("S3").Value = 0.1807312169294119
$ws.Range("T3").Value = 0.1807312169294119
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 37.63904266666666
$ws.Range("H4").Value = 112.917128
$ws.Range("I4").Value = 0.4850220755088102
$ws.Range("J4").Value = 0.4850220755088102
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.36964133333333
$ws.Range("N4").Value = 76.108924
$ws.Range("O4").Value = 0.2337290805561598
$ws.Range("P4").Value = 0.2337290805561598
$ws.Range("Q4").Value = 954.8890125833634
$ws.Range("R4").Value = 8594.001113250271
$ws.Range("S4").Value = 0.1133637637581145
$ws.Range("T4").Value = 0.1133637637581145
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 37.63904266666666
$ws.Range("H5").Value = 112.917128
$ws.Range("I5").Value = 0.4850220755088102
$ws.Range("J5").Value = 0.4850220755088102
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.49855033333333
$ws.Range("N5").Value = 34.495651
$ws.Range("O5").Value = 0.1059354983315251
$ws.Range("P5").Value = 0.1059354983315251
$ws.Range("Q5").Value = 432.7944266011475
$ws.Range("R5").Value = 3895.149839410328
$ws.Range("S5").Value = 0.05138105527081642
$ws.Range("T5").Value = 0.05138105527081643
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.57434666666667
$ws.Range("H6").Value = 52.72304
$ws.Range("I6").Value = 0.2264655392929762
$ws.Range("J6").Value = 0.2264655392929762
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.22896466666667
$ws.Range("N6").Value = 93.686894
$ws.Range("O6").Value = 0.2877106972998646
$ws.Range("P6").Value = 0.2877106972998646
$ws.Range("Q6").Value = 548.8286510930844
$ws.Range("R6").Value = 4939.457859837759
$ws.Range("S6").Value = 0.06515655822437207
$ws.Range("T6").Value = 0.06515655822437207
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.57434666666667
$ws.Range("H7").Value = 52.72304
$ws.Range("I7").Value = 0.2264655392929762
$ws.Range("J7").Value = 0.2264655392929762
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 40.44578266666667
$ws.Range("N7").Value = 121.337348
$ws.Range("O7").Value = 0.3726247238124506
$ws.Range("P7").Value = 0.3726247238124505
$ws.Range("Q7").Value = 710.8082057886578
$ws.Range("R7").Value = 6397.273852097919
$ws.Range("S7").Value = 0.08438665903208294
$ws.Range("T7").Value = 0.08438665903208292
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.57434666666667
$ws.Range("H8").Value = 52.72304
$ws.Range("I8").Value = 0.2264655392929762
$ws.Range("J8").Value = 0.2264655392929762
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.36964133333333
$ws.Range("N8").Value = 76.108924
$ws.Range("O8").Value = 0.2337290805561598
$ws.Range("P8").Value = 0.2337290805561598
$ws.Range("Q8").Value = 445.8548716009956
$ws.Range("R8").Value = 4012.69384440896
$ws.Range("S8").Value = 0.0529315822766022
$ws.Range("T8").Value = 0.0529315822766022
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.57434666666667
$ws.Range("H9").Value = 52.72304
$ws.Range("I9").Value = 0.2264655392929762
$ws.Range("J9").Value = 0.2264655392929762
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.49855033333333
$ws.Range("N9").Value = 34.495651
$ws.Range("O9").Value = 0.1059354983315251
$ws.Range("P9").Value = 0.1059354983315251
$ws.Range("Q9").Value = 202.0795097221156
$ws.Range("R9").Value = 1818.71558749904
$ws.Range("S9").Value = 0.02399073975991902
$ws.Range("T9").Value = 0.02399073975991902
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.168158
$ws.Range("H10").Value = 57.504474
$ws.Range("I10").Value = 0.2470036195972184
$ws.Range("J10").Value = 0.2470036195972184
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.22896466666667
$ws.Range("N10").Value = 93.686894
$ws.Range("O10").Value = 0.2877106972998646
$ws.Range("P10").Value = 0.2877106972998646
$ws.Range("Q10").Value = 598.6017289070841
$ws.Range("R10").Value = 5387.415560163756
$ws.Range("S10").Value = 0.0710655836299062
$ws.Range("T10").Value = 0.0710655836299062
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 19.168158
$ws.Range("H11").Value = 57.504474
$ws.Range("I11").Value = 0.2470036195972184
$ws.Range("J11").Value = 0.2470036195972184
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 40.44578266666667
$ws.Range("N11").Value = 121.337348
$ws.Range("O11").Value = 0.3726247238124506
$ws.Range("P11").Value = 0.3726247238124505
$ws.Range("Q11").Value = 775.2711525883281
$ws.Range("R11").Value = 6977.440373294951
$ws.Range("S11").Value = 0.09203965553308911
$ws.Range("T11").Value = 0.0920396555330891
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 19.168158
$ws.Range("H12").Value = 57.504474
$ws.Range("I12").Value = 0.2470036195972184
$ws.Range("J12").Value = 0.2470036195972184
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 25.36964133333333
$ws.Range("N12").Value = 76.108924
$ws.Range("O12").Value = 0.2337290805561598
$ws.Range("P12").Value = 0.2337290805561598
$ws.Range("Q12").Value = 486.2892934806641
$ws.Range("R12").Value = 4376.603641325976
$ws.Range("S12").Value = 0.0577319289025013
$ws.Range("T12").Value = 0.0577319289025013
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 19.168158
$ws.Range("H13").Value = 57.504474
$ws.Range("I13").Value = 0.2470036195972184
$ws.Range("J13").Value = 0.2470036195972184
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 11.49855033333333
$ws.Range("N13").Value = 34.495651
$ws.Range("O13").Value = 0.1059354983315251
$ws.Range("P13").Value = 0.1059354983315251
$ws.Range("Q13").Value = 220.406029560286
$ws.Range("R13").Value = 1983.654266042574
$ws.Range("S13").Value = 0.0261664515317218
$ws.Range("T13").Value = 0.0261664515317218
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.221194
$ws.Range("H14").Value = 9.663582
$ws.Range("I14").Value = 0.04150876560099527
$ws.Range("J14").Value = 0.04150876560099527
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.22896466666667
$ws.Range("N14").Value = 93.686894
$ws.Range("O14").Value = 0.2877106972998646
$ws.Range("P14").Value = 0.2877106972998646
$ws.Range("Q14").Value = 100.5945536104787
$ws.Range("R14").Value = 905.3509824943079
$ws.Range("S14").Value = 0.01194251589511898
$ws.Range("T14").Value = 0.01194251589511898
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.221194
$ws.Range("H15").Value = 9.663582
$ws.Range("I15").Value = 0.04150876560099527
$ws.Range("J15").Value = 0.04150876560099527
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 40.44578266666667
$ws.Range("N15").Value = 121.337348
$ws.Range("O15").Value = 0.3726247238124506
$ws.Range("P15").Value = 0.3726247238124505
$ws.Range("Q15").Value = 130.2837124511707
$ws.Range("R15").Value = 1172.553412060536
$ws.Range("S15").Value = 0.01546719231786661
$ws.Range("T15").Value = 0.01546719231786661
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.221194
$ws.Range("H16").Value = 9.663582
$ws.Range("I16").Value = 0.04150876560099527
$ws.Range("J16").Value = 0.04150876560099527
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 25.36964133333333
$ws.Range("N16").Value = 76.108924
$ws.Range("O16").Value = 0.2337290805561598
$ws.Range("P16").Value = 0.2337290805561598
$ws.Range("Q16").Value = 81.72053644508534
$ws.Range("R16").Value = 735.484828005768
$ws.Range("S16").Value = 0.009701805618941777
$ws.Range("T16").Value = 0.009701805618941777
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.221194
$ws.Range("H17").Value = 9.663582
$ws.Range("I17").Value = 0.04150876560099527
$ws.Range("J17").Value = 0.04150876560099527
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 11.49855033333333
$ws.Range("N17").Value = 34.495651
$ws.Range("O17").Value = 0.1059354983315251
$ws.Range("P17").Value = 0.1059354983315251
$ws.Range("Q17").Value = 37.03906134243134
$ws.Range("R17").Value = 333.351552081882
$ws.Range("S17").Value = 0.004397251769067903
$ws.Range("T17").Value = 0.004397251769067903
